# Applies the updated crypto price/volume figures to Sheet1 (rows 2-50).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.034.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.501.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.56%  "
$ws.Range("E7").Value = "  -0.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.498.58"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.191"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.579"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.064.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "608.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.502.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.145.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.873"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -14.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.93%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  -1.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "643.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +16.05%  "
$ws.Range("E34").Value = "  -3.73%  "
$ws.Range("E36").Value = "  +3.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0987"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0468"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.351.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0736"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.82%  "
$ws.Range("E45").Value = "  -4.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "32.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("E48").Value = "  -2.21%  "
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.31%  "
